$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before "disponible" (D) for the new "localidad" field.
# Excel copies the formatting of the column to the left (direccion) into the
# freshly inserted column, which also shifts disponible/plazas/observaciones
# one column to the right.
$ws.Range("D1").EntireColumn.Insert()
$ws.Range("D1").Value = "localidad"

# The address used to include the locality ("...,  La Matanza"); split it so
# the locality lives in its own column and the street address is shorter.
$ws.Range("C2").Value = "El Pampero 5790"
$ws.Range("D2").Value = "la matanza"

# "observaciones" is still the last header, now shifted out to column G.
$ws.Range("G1").Value = "observaciones"

# The row is now shorter, since the wrapped address text is shorter than
# before - let the row height shrink to fit the new (shorter) content.
$ws.Rows("2:2").RowHeight = 23.85

# Add the new driver row.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "jorge"
$ws.Range("C3").Value = "tomas valle 5890"
$ws.Range("D3").Value = "la matanza"
$ws.Range("E3").Value = "si"
$ws.Range("F3").Value = 4

# Match the font used in the wrapped address column for the (still empty)
# observaciones cell of the new row, without wrapping.
$ws.Range("G3").Font.Name = "Arial"
$ws.Range("G3").Font.Size = 10
$ws.Range("G3").Font.Underline = 2
$ws.Range("G3").WrapText = $false

$ws.Range("G3").Select()
